$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# --- Replace the text of the first paragraph, merging/removing the trailing
#     space-only run in the process -----------------------------------------
$fullRange = $p1.Range
$textRange = $d.Range($fullRange.Start, $fullRange.End - 1)
$textRange.Text = "**ID__AFFARS_SUBPART_5302_1__ID**"

# --- Paragraph formatting: add a paragraph border (space-only, no line) and
#     widen the left indent --------------------------------------------------
$p1 = $d.Paragraphs(1)
$pRange = $p1.Range
$borders = $pRange.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

$p1.Format.LeftIndent = 11.25
